$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between row 16 and row 17
$ws.Range("E16").Value = "1907"
$ws.Range("E17").Value = "1906"

# Swap the "Valor Mora" values between row 16 and row 17
$ws.Range("F16").Value = 48000
$ws.Range("F17").Value = 33600
